$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Price column (D) holds numeric-looking strings that must stay text (inline/shared
# strings in the source file), matching the source data's type. Temporarily mark the
# whole column as Text before writing, then clear the formatting residue afterwards so
# the cell style indices end up identical to the original (style 0 / no explicit style).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.974.21'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').Value = '1.597.93'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '210.80'
$ws.Range('E5').Value = '  +2.13%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('D8').Value = '0.245'
$ws.Range('E8').Value = '  +0.97%  '
$ws.Range('D9').Value = '0.0611'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').Value = '18.05'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('D12').Value = '1.822.52'
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('D13').Value = '1.604.82'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').Value = '3.99'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = '0.513'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = '25.978.03'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').Value = '60.11'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '0.0₃0721'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D20').Value = '201.10'
$ws.Range('E20').Value = '  +7.61%  '
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('D22').Value = '9.25'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '6.00'
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('E24').Value = '  +7.36%  '
$ws.Range('D25').Value = '141.56'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -7.37%  '
$ws.Range('D28').Value = '15.10'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('E29').Value = '  +1.03%  '
$ws.Range('D30').Value = '1.16'
$ws.Range('E30').Value = '  +1.91%  '
$ws.Range('D31').Value = '0.0475'
$ws.Range('E31').Value = '  +2.11%  '
$ws.Range('D32').Value = '3.09'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').Value = '2.95'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').Value = '1.47'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('D36').Value = '1.121.79'
$ws.Range('E36').Value = '  +2.97%  '
$ws.Range('E37').Value = '  +9.90%  '
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = '0.786'
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('D41').Value = '0.489'
$ws.Range('E41').Value = '  -1.55%  '
$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.733.89'
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.12'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').Value = '92.67'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').Value = '53.34'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = '7.19'
$ws.Range('E51').Value = '  -0.22%  '

$ws.Range("D2:D51").ClearFormats()

